$wb = $excel.ActiveWorkbook

# Sheet "建物" (building): row 2, column I (property_category) "land" -> "building"
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"

# Sheet "汽車" (car): row 2, column H (property_category) "land" -> "car"
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
